# BulunTe_IphoneMonitor.xlsx - append 6 new daily-tracking rows (29-34)
# covering 1/27/24 .. 2/1/24, mirroring the existing layout of rows 18-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cellAddr, $sourceAddr, $text) {
    # Clone the number format + alignment of an existing date cell (style index 2)
    # without ever letting Excel "helpfully" re-interpret the literal text as a
    # real date serial number (which would create a brand-new style).
    $ws.Range($sourceAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($cellAddr).Formula = "=""" + $text + """"
    $ws.Range($cellAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)   # xlPasteValues (collapses formula -> literal text)
}

# ---- Row 29 : 1/27/24 ----
Set-DateText "A29" "A18" "1/27/24"
$ws.Range("B29").Value = "59min"
$ws.Range("D29").Value = "22min"
$ws.Range("F29").Value = 34
$ws.Range("G29").NumberFormat = "h:mm"
$ws.Range("G29").Value = 0.05486111111111111

# ---- Row 30 : 1/28/24 ----
Set-DateText "A30" "A18" "1/28/24"
$ws.Range("B30").Value = "26min"
$ws.Range("D30").Value = "15min"
$ws.Range("F30").Value = 31
$ws.Range("G30").NumberFormat = "h:mm"
$ws.Range("G30").Value = 0.0020833333333333333

# ---- Row 31 : 1/29/24 ----
Set-DateText "A31" "A18" "1/29/24"
$ws.Range("B31").Value = "2h12min"
$ws.Range("D31").Value = "23min"
$ws.Range("F31").Value = 76
$ws.Range("G31").NumberFormat = "h:mm"
$ws.Range("G31").Value = 0.32222222222222224

# ---- Row 32 : 1/30/24 ----
Set-DateText "A32" "A18" "1/30/24"
$ws.Range("B32").Value = "1h34min"
$ws.Range("D32").Value = "37min"
$ws.Range("F32").Value = 85
$ws.Range("G32").NumberFormat = "h:mm"
$ws.Range("G32").Value = 0.3201388888888889

# ---- Row 33 : 1/31/24 ----
Set-DateText "A33" "A18" "1/31/24"
$ws.Range("B33").Value = "53min"
$ws.Range("D33").Value = "38min"
$ws.Range("F33").Value = 57
$ws.Range("G33").NumberFormat = "h:mm"
$ws.Range("G33").Value = 0.02152777777777778

# ---- Row 34 : 2/1/24 ----
Set-DateText "A34" "A18" "2/1/24"
$ws.Range("B34").Value = "1h10min"
$ws.Range("D34").Value = "24min"
$ws.Range("F34").Value = 84
$ws.Range("G34").NumberFormat = "h:mm"
$ws.Range("G34").Value = 0.009722222222222222

# ---- C & E columns: extend the "minutes parsed from Bxx/Dxx" shared formulas ----
$cFormula = "=IF(ISERROR(FIND(""h"", B29)), 0, LEFT(B29, FIND(""h"", B29)-1)*60) + IF(ISERROR(FIND(""min"", B29)), 0, MID(B29, IF(ISERROR(FIND(""h"", B29)), 1, FIND(""h"", B29)+1), FIND(""min"", B29) - IF(ISERROR(FIND(""h"", B29)), 1, FIND(""h"", B29)+1)))"
$ws.Range("C29:C34").Formula = $cFormula

$eFormula = "=IF(ISERROR(FIND(""h"", D29)), 0, LEFT(D29, FIND(""h"", D29)-1)*60) + IF(ISERROR(FIND(""min"", D29)), 0, MID(D29, IF(ISERROR(FIND(""h"", D29)), 1, FIND(""h"", D29)+1), FIND(""min"", D29) - IF(ISERROR(FIND(""h"", D29)), 1, FIND(""h"", D29)+1)))"
$ws.Range("E29:E34").Formula = $eFormula

# ---- view state: scroll position + active selection ----
$excel.Goto($ws.Range("A7"), $false)
$ws.Range("G35").Select()
